# Insert a new data row at row 239 (pushing the existing rows 239-296 down
# to 240-297) and populate it with a new price observation record.
# This mirrors the source diff, where the sheet dimension grows from
# A1:R296 to A1:R297 and a "new" weekly reading is spliced in above the
# previously-first record for that date range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 239..296 down to 240..297, leaving row 239 empty.
$ws.Rows.Item(239).Insert()

# Populate the newly inserted row 239 with the new record.
$ws.Range("A239").Value = 10
$ws.Range("B239").Value = "Vega Modelo de Temuco"
$ws.Range("C239").Value = "La Araucanía"
$ws.Range("D239").Value = 44642
$ws.Range("E239").Value = 9
$ws.Range("F239").Value = 100112009
$ws.Range("G239").Value = "Acelga"
$ws.Range("H239").Value = "Sin especificar"
$ws.Range("I239").Value = "Primera"
$ws.Range("J239").Value = 45
$ws.Range("K239").Value = 8000
$ws.Range("L239").Value = 8000
$ws.Range("M239").Value = 8000
$ws.Range("N239").Value = "`$/docena de atados (12 kilos)"
$ws.Range("O239").Value = "Provincia de Cautín"
$ws.Range("P239").Value = 667
$ws.Range("Q239").Value = 12
$ws.Range("R239").Value = "Hortaliza"
